# Update cryptocurrency price (D) and volume-change (E) columns
# with freshly scraped values, preserving original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value = "37.781.87"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "

$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value = "2.076.64"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("E4").Value = "  +0.01%  "

$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value = "232.79"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "

$ws.Range("E6").Value = "  -0.23%  "

$cD = $ws.Range("D7")
$cD.NumberFormat = "@"
$cD.Value = "58.62"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.26%  "

$cD = $ws.Range("D10")
$cD.NumberFormat = "@"
$cD.Value = "0.0785"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E10").Value = "  -1.16%  "

$ws.Range("E11").Value = "  +3.03%  "

$cD = $ws.Range("D12")
$cD.NumberFormat = "@"
$cD.Value = "14.86"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "

$cD = $ws.Range("D13")
$cD.NumberFormat = "@"
$cD.Value = "2.382.69"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E13").Value = "  -0.55%  "

$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value = "21.08"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E14").Value = "  -1.12%  "

$cD = $ws.Range("D15")
$cD.NumberFormat = "@"
$cD.Value = "0.784"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("E16").Value = "  +0.69%  "

$cD = $ws.Range("D17")
$cD.NumberFormat = "@"
$cD.Value = "2.097.63"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E17").Value = "  +0.53%  "

$cD = $ws.Range("D18")
$cD.NumberFormat = "@"
$cD.Value = "37.691.46"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E18").Value = "  -0.13%  "

$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value = "6.13"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E19").Value = "  -2.74%  "

$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value = "71.42"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "

$cD = $ws.Range("D21")
$cD.NumberFormat = "@"
$cD.Value = "0.0₃0840"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E21").Value = "  +1.17%  "

$cD = $ws.Range("D22")
$cD.NumberFormat = "@"
$cD.Value = "229.21"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  -0.88%  "

$cD = $ws.Range("D25")
$cD.NumberFormat = "@"
$cD.Value = "2.39"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "

$cD = $ws.Range("D26")
$cD.NumberFormat = "@"
$cD.Value = "9.69"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E26").Value = "  +6.85%  "

$cD = $ws.Range("D27")
$cD.NumberFormat = "@"
$cD.Value = "171.94"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("E28").Value = "  -1.22%  "

$cD = $ws.Range("D29")
$cD.NumberFormat = "@"
$cD.Value = "19.45"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("E30").Value = "  -2.19%  "

$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("E32").Value = "  +0.01%  "

$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value = "0.0632"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E33").Value = "  -0.14%  "

$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value = "4.67"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "

$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("E36").Value = "  -0.84%  "

$ws.Range("E37").Value = "  -2.56%  "

$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value = "0.999"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "

$cD = $ws.Range("D39")
$cD.NumberFormat = "@"
$cD.Value = "5.41"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E39").Value = "  -0.45%  "

$cD = $ws.Range("D40")
$cD.NumberFormat = "@"
$cD.Value = "0.0233"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E40").Value = "  +7.97%  "

$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value = "100.97"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E41").Value = "  +1.85%  "

$ws.Range("E42").Value = "  -1.26%  "

$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("E44").Value = "  +4.94%  "

$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value = "1.445.86"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("E47").Value = "  -1.46%  "

$ws.Range("E48").Value = "  -4.99%  "

$ws.Range("E49").Value = "  -1.89%  "

$ws.Range("E50").Value = "  -1.95%  "

$cD = $ws.Range("D51")
$cD.NumberFormat = "@"
$cD.Value = "2.268.34"
$cD.NumberFormat = "General"
$cD.Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
